$d = $word.ActiveDocument

# 1) Update the date text 09/09/2021 -> 09/10/2021
$d.Content.Find.Execute("09/09/2021", $true, $false, $false, $false, $false, $true, 1, $false, "09/10/2021", 2)

# 2) Remove the old mid-document _GoBack bookmark
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
